$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update individual F-column values (imputed error values)
$ws.Range("F3").Value = 17.64
$ws.Range("F5").Value = ""
$ws.Range("F21").Value = 16.58
$ws.Range("F23").Value = ""

# Remove the "RM 232" row (originally row 26); rows below shift up
$ws.Rows("26").Delete()

# After the shift, "SC 92" (originally row 28) is now row 27; remove it too
$ws.Rows("27").Delete()

# The row that is now 32 (previously "SC 193" at row 34) gets its F value filled in
$ws.Range("F32").Value = 17.39
